$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (player, position, team) for rows 2..18 (A2:C18)
$data = @(
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Luke Kennard", "SG", "Memphis Grizzlies"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# The old sheet had 18 data rows (19 total incl. header); the new sheet only has
# 17 data rows (18 total), so remove the now-superfluous last row entirely.
$ws.Rows.Item(19).Delete()
